# Apply the TestPlan update: mark all existing test rows (3-24) as executed
# with Actual result / Notes / Status (PASS, highlighted green).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

for ($r = 3; $r -le 24; $r++) {
    # Column order matters for shared-string insertion order: Notes ("Done on
    # 15.04.19") must be added to the shared string table before Actual result
    # ("As expected").
    $ws.Range("G$r").Value = "Done on 15.04.19"
    $ws.Range("F$r").Value = "As expected"
    $ws.Range("H$r").Value = "PASS"
}

# Highlight the Status column (H3:H24) with a green fill and the same
# top/wrap alignment used elsewhere in the table.
$statusRange = $ws.Range("H3:H24")
$statusRange.Interior.Color = 5287936
$statusRange.WrapText = $true
$statusRange.VerticalAlignment = -4160

# Restore the saved cursor/selection position.
[void]$ws.Range("G14").Select()
